$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Result" header column - match the existing header style (bold / fill / border)
# used by the other header cells (A1, B1) before writing the text.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Result"

# Per-row PASS/FAIL results
$ws.Range("C2").Value = "PASS"
$ws.Range("C3").Value = "PASS"
$ws.Range("C4").Value = "FAIL"
$ws.Range("C5").Value = "PASS"
$ws.Range("C6").Value = "PASS"

# Update selection to mirror the author's saved selection state
$ws.Range("C2:F11").Select()
